$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix typo: "Sand-Attack" -> "Sand Attack" in the Galarian Zigzagoon move list ---
$ws.Range("E106").Value = "Snarl, Headbutt, Sand Attack, Leer"

# --- 2. Add new trainer block: TRAINER_MARY (Meditite, Cubone), then move END marker down ---
# A style template cell (matches the "s=1" styling used on header / first-pokemon rows
# throughout the sheet, e.g. A1 / A2:E2).
$styleTemplate = $ws.Range("A1")

# Row 109: new trainer name header "TRAINER_MARY" (default/unstyled, like row 104 "TRAINER_NORMAN_1")
$ws.Range("A109").Value = "TRAINER_MARY"

# Row 110: column headers for the new trainer block
$ws.Range("A110").Value = "species"
$ws.Range("B110").Value = "lvl"
$ws.Range("C110").Value = "iv"
$ws.Range("D110").Value = "heldItem"
$ws.Range("E110").Value = "moves"
$ws.Range("A110:E110").Style = $styleTemplate.Style

# Row 111: first pokemon - Meditite, lvl 5, blank iv cell present
$ws.Range("A111").Value = "Meditite"
$ws.Range("B111").Value = 5
$ws.Range("C111").Style = $styleTemplate.Style
$ws.Range("A111:B111").Style = $styleTemplate.Style

# Row 112: second pokemon - Cubone, lvl 5 (default style, like other trailing pokemon rows)
$ws.Range("A112").Value = "Cubone"
$ws.Range("B112").Value = 5

# Row 114: END marker moved down from its old position (row 111) to row 114
$ws.Range("A114").Value = "END"
$ws.Range("A114").Style = $styleTemplate.Style

# --- 3. Update the active selection to reflect where editing left off ---
$ws.Range("C112").Select()
